# "fix some minor stuff" - tweak a handful of ergonomics/strength inputs on
# the m4-barrels sheet and update the N-column computed totals accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("m4-barrels")

# Row 5 (HK MP5 Wide Tropical): ergonomics 0.2 -> 0.18
$ws.Range("D5").Value = 0.18

# Row 6 (HK MP5SD Polymer Handguard): strength 10 -> 11, ergonomics 0.24 -> 0.22
$ws.Range("C6").Value = 11
$ws.Range("D6").Value = 0.22

# Row 7 (HK MP5 Slim Checkered Handguard): ergonomics 0.18 -> 0.16
$ws.Range("D7").Value = 0.16

# Row 8 (Surefire 628LMF-B MP5 Handguard): ergonomics 0.23 -> 0.21
$ws.Range("D8").Value = 0.21

# Row 9 (Midwest Industries MP5 MLOK Handguard): ergonomics 0.17 -> 0.15
$ws.Range("D9").Value = 0.15

# The N column holds a shared formula (C-D*20-E*0.8-F*0.6-H*5+I*10+J/300) that
# recalculates automatically from the edits above.

# Leave the cursor where the author left it when they saved.
$ws.Range("J6").Select()
